$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Col4a1"
$ws.Cells.Item(2,3).Value = "Itgb8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 216.5832213333333
$ws.Cells.Item(2,8).Value = 649.749664
$ws.Cells.Item(2,9).Value = 0.4331411212367192
$ws.Cells.Item(2,10).Value = 0.4331411212367192
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.027767
$ws.Cells.Item(2,14).Value = 0.083301
$ws.Cells.Item(2,15).Value = 0.002463719941166009
$ws.Cells.Item(2,16).Value = 0.002463719941166009
$ws.Cells.Item(2,17).Value = 6.013866306762667
$ws.Cells.Item(2,18).Value = 54.124796760864
$ws.Cells.Item(2,19).Value = 0.001067138417729909
$ws.Cells.Item(2,20).Value = 0.001067138417729909

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Col4a1"
$ws.Cells.Item(3,3).Value = "Itgb8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 216.5832213333333
$ws.Cells.Item(3,8).Value = 649.749664
$ws.Cells.Item(3,9).Value = 0.4331411212367192
$ws.Cells.Item(3,10).Value = 0.4331411212367192
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.987076
$ws.Cells.Item(3,14).Value = 11.961228
$ws.Cells.Item(3,15).Value = 0.3537666527944829
$ws.Cells.Item(3,16).Value = 0.3537666527944829
$ws.Cells.Item(3,17).Value = 863.5337637808213
$ws.Cells.Item(3,18).Value = 7771.803874027392
$ws.Cells.Item(3,19).Value = 0.1532308846475635
$ws.Cells.Item(3,20).Value = 0.1532308846475635

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Col4a1"
$ws.Cells.Item(4,3).Value = "Itgb8"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 216.5832213333333
$ws.Cells.Item(4,8).Value = 649.749664
$ws.Cells.Item(4,9).Value = 0.4331411212367192
$ws.Cells.Item(4,10).Value = 0.4331411212367192
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 7.255512666666667
$ws.Cells.Item(4,14).Value = 21.766538
$ws.Cells.Item(4,15).Value = 0.643769627264351
$ws.Cells.Item(4,16).Value = 0.643769627264351
$ws.Cells.Item(4,17).Value = 1571.42230577147
$ws.Cells.Item(4,18).Value = 14142.80075194323
$ws.Cells.Item(4,19).Value = 0.2788430981714258
$ws.Cells.Item(4,20).Value = 0.2788430981714258

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Col4a1"
$ws.Cells.Item(5,3).Value = "Itgb8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 161.954974
$ws.Cells.Item(5,8).Value = 485.864922
$ws.Cells.Item(5,9).Value = 0.3238910133313607
$ws.Cells.Item(5,10).Value = 0.3238910133313606
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.027767
$ws.Cells.Item(5,14).Value = 0.083301
$ws.Cells.Item(5,15).Value = 0.002463719941166009
$ws.Cells.Item(5,16).Value = 0.002463719941166009
$ws.Cells.Item(5,17).Value = 4.497003763058
$ws.Cells.Item(5,18).Value = 40.473033867522
$ws.Cells.Item(5,19).Value = 0.0007979767483089391
$ws.Cells.Item(5,20).Value = 0.000797976748308939

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Col4a1"
$ws.Cells.Item(6,3).Value = "Itgb8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 161.954974
$ws.Cells.Item(6,8).Value = 485.864922
$ws.Cells.Item(6,9).Value = 0.3238910133313607
$ws.Cells.Item(6,10).Value = 0.3238910133313606
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.987076
$ws.Cells.Item(6,14).Value = 11.961228
$ws.Cells.Item(6,15).Value = 0.3537666527944829
$ws.Cells.Item(6,16).Value = 0.3537666527944829
$ws.Cells.Item(6,17).Value = 645.726789916024
$ws.Cells.Item(6,18).Value = 5811.541109244215
$ws.Cells.Item(6,19).Value = 0.1145818396564487
$ws.Cells.Item(6,20).Value = 0.1145818396564487

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Col4a1"
$ws.Cells.Item(7,3).Value = "Itgb8"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 161.954974
$ws.Cells.Item(7,8).Value = 485.864922
$ws.Cells.Item(7,9).Value = 0.3238910133313607
$ws.Cells.Item(7,10).Value = 0.3238910133313606
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 7.255512666666667
$ws.Cells.Item(7,14).Value = 21.766538
$ws.Cells.Item(7,15).Value = 0.643769627264351
$ws.Cells.Item(7,16).Value = 0.643769627264351
$ws.Cells.Item(7,17).Value = 1175.066365286671
$ws.Cells.Item(7,18).Value = 10575.59728758004
$ws.Cells.Item(7,19).Value = 0.208511196926603
$ws.Cells.Item(7,20).Value = 0.208511196926603

$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Col4a1"
$ws.Cells.Item(8,3).Value = "Itgb8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.4608033333333333
$ws.Cells.Item(8,8).Value = 1.38241
$ws.Cells.Item(8,9).Value = 0.0009215527926904059
$ws.Cells.Item(8,10).Value = 0.0009215527926904059
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.027767
$ws.Cells.Item(8,14).Value = 0.083301
$ws.Cells.Item(8,15).Value = 0.002463719941166009
$ws.Cells.Item(8,16).Value = 0.002463719941166009
$ws.Cells.Item(8,17).Value = 0.01279512615666667
$ws.Cells.Item(8,18).Value = 0.11515613541
$ws.Cells.Item(8,19).Value = 0.000002270447992188578
$ws.Cells.Item(8,20).Value = 0.000002270447992188578

$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Col4a1"
$ws.Cells.Item(9,3).Value = "Itgb8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.4608033333333333
$ws.Cells.Item(9,8).Value = 1.38241
$ws.Cells.Item(9,9).Value = 0.0009215527926904059
$ws.Cells.Item(9,10).Value = 0.0009215527926904059
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.987076
$ws.Cells.Item(9,14).Value = 11.961228
$ws.Cells.Item(9,15).Value = 0.3537666527944829
$ws.Cells.Item(9,16).Value = 0.3537666527944829
$ws.Cells.Item(9,17).Value = 1.837257911053333
$ws.Cells.Item(9,18).Value = 16.53532119948
$ws.Cells.Item(9,19).Value = 0.000326014646843493
$ws.Cells.Item(9,20).Value = 0.000326014646843493

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Col4a1"
$ws.Cells.Item(10,3).Value = "Itgb8"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.4608033333333333
$ws.Cells.Item(10,8).Value = 1.38241
$ws.Cells.Item(10,9).Value = 0.0009215527926904059
$ws.Cells.Item(10,10).Value = 0.0009215527926904059
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 7.255512666666667
$ws.Cells.Item(10,14).Value = 21.766538
$ws.Cells.Item(10,15).Value = 0.643769627264351
$ws.Cells.Item(10,16).Value = 0.643769627264351
$ws.Cells.Item(10,17).Value = 3.343364421842222
$ws.Cells.Item(10,18).Value = 30.09027979658
$ws.Cells.Item(10,19).Value = 0.0005932676978547244
$ws.Cells.Item(10,20).Value = 0.0005932676978547244

$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Col4a1"
$ws.Cells.Item(11,3).Value = "Itgb8"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 121.0302313333333
$ws.Cells.Item(11,8).Value = 363.090694
$ws.Cells.Item(11,9).Value = 0.2420463126392298
$ws.Cells.Item(11,10).Value = 0.2420463126392298
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.027767
$ws.Cells.Item(11,14).Value = 0.083301
$ws.Cells.Item(11,15).Value = 0.002463719941166009
$ws.Cells.Item(11,16).Value = 0.002463719941166009
$ws.Cells.Item(11,17).Value = 3.360646433432667
$ws.Cells.Item(11,18).Value = 30.245817900894
$ws.Cells.Item(11,19).Value = 0.0005963343271349728
$ws.Cells.Item(11,20).Value = 0.0005963343271349727

$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Col4a1"
$ws.Cells.Item(12,3).Value = "Itgb8"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 121.0302313333333
$ws.Cells.Item(12,8).Value = 363.090694
$ws.Cells.Item(12,9).Value = 0.2420463126392298
$ws.Cells.Item(12,10).Value = 0.2420463126392298
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 3.987076
$ws.Cells.Item(12,14).Value = 11.961228
$ws.Cells.Item(12,15).Value = 0.3537666527944829
$ws.Cells.Item(12,16).Value = 0.3537666527944829
$ws.Cells.Item(12,17).Value = 482.5567306235814
$ws.Cells.Item(12,18).Value = 4343.010575612232
$ws.Cells.Item(12,19).Value = 0.08562791384362728
$ws.Cells.Item(12,20).Value = 0.08562791384362727

$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Col4a1"
$ws.Cells.Item(13,3).Value = "Itgb8"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 121.0302313333333
$ws.Cells.Item(13,8).Value = 363.090694
$ws.Cells.Item(13,9).Value = 0.2420463126392298
$ws.Cells.Item(13,10).Value = 0.2420463126392298
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 7.255512666666667
$ws.Cells.Item(13,14).Value = 21.766538
$ws.Cells.Item(13,15).Value = 0.643769627264351
$ws.Cells.Item(13,16).Value = 0.643769627264351
$ws.Cells.Item(13,17).Value = 878.136376488597
$ws.Cells.Item(13,18).Value = 7903.227388397372
$ws.Cells.Item(13,19).Value = 0.1558220644684676
$ws.Cells.Item(13,20).Value = 0.1558220644684675
